$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D cells keep their original text formatting so numeric-looking
# strings like "598.98" are not converted to actual numbers.
$ws.Range('D2:D51').NumberFormat = '@'

$ws.Range('D2').Value = '64.888.46'
$ws.Range('E2').Value = '  -0.43%  '
$ws.Range('D3').Value = '3.564.78'
$ws.Range('E3').Value = '  +2.54%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '598.98'
$ws.Range('E5').Value = '  +1.67%  '
$ws.Range('D6').Value = '135.32'
$ws.Range('E6').Value = '  -1.10%  '
$ws.Range('D7').Value = '3.562.91'
$ws.Range('E7').Value = '  +2.52%  '
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('E9').Value = '  +0.58%  '
$ws.Range('E10').Value = '  +0.37%  '
$ws.Range('D11').Value = '6.95'
$ws.Range('E11').Value = '  -2.63%  '
$ws.Range('E12').Value = '  +0.32%  '
$ws.Range('D13').Value = '4.164.65'
$ws.Range('E13').Value = '  +2.42%  '
$ws.Range('E14').Value = '  -0.25%  '
$ws.Range('D15').Value = '3.561.92'
$ws.Range('E15').Value = '  +2.53%  '
$ws.Range('D16').Value = '27.01'
$ws.Range('E16').Value = '  +1.60%  '
$ws.Range('E17').Value = '  +0.48%  '
$ws.Range('D18').Value = '64.482.57'
$ws.Range('E18').Value = '  -1.00%  '
$ws.Range('D19').Value = '10.00'
$ws.Range('E19').Value = '  +3.13%  '
$ws.Range('D20').Value = '14.34'
$ws.Range('E20').Value = '  +3.23%  '
$ws.Range('E21').Value = '  +0.90%  '
$ws.Range('D22').Value = '388.52'
$ws.Range('E22').Value = '  -0.11%  '
$ws.Range('E23').Value = '  +4.26%  '
$ws.Range('D24').Value = '3.704.63'
$ws.Range('E24').Value = '  +2.37%  '
$ws.Range('D25').Value = '73.97'
$ws.Range('E25').Value = '  +1.88%  '
$ws.Range('D27').Value = '0.0000114'
$ws.Range('E27').Value = '  +4.00%  '
$ws.Range('D28').Value = '7.68'
$ws.Range('E28').Value = '  +4.92%  '
$ws.Range('E29').Value = '  +0.00%  '
$ws.Range('D30').Value = '2.28'
$ws.Range('E30').Value = '  +2.98%  '
$ws.Range('D31').Value = '8.41'
$ws.Range('E31').Value = '  +2.89%  '
$ws.Range('D32').Value = '1.47'
$ws.Range('E32').Value = '  +23.31%  '
$ws.Range('D33').Value = '3.561.17'
$ws.Range('E33').Value = '  +1.88%  '
$ws.Range('D34').Value = '24.00'
$ws.Range('E34').Value = '  +3.68%  '
$ws.Range('E35').Value = '  +0.00%  '
$ws.Range('D36').Value = '0.143'
$ws.Range('E36').Value = '  +0.41%  '
$ws.Range('B37').Value = 'Aptos'
$ws.Range('C37').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D37').Value = '6.91'
$ws.Range('E37').Value = '  +1.42%  '
$ws.Range('B38').Value = 'Monero'
$ws.Range('C38').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D38').Value = '168.61'
$ws.Range('E38').Value = '  -1.16%  '
$ws.Range('E39').Value = '  +4.62%  '
$ws.Range('E40').Value = '  +5.64%  '
$ws.Range('D41').Value = '0.0804'
$ws.Range('E41').Value = '  +3.31%  '
$ws.Range('B42').Value = 'Mantle'
$ws.Range('C42').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D42').Value = '0.824'
$ws.Range('E42').Value = '  +1.57%  '
$ws.Range('B43').Value = 'EnergySwap'
$ws.Range('C43').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D43').Value = '26.74'
$ws.Range('E43').Value = '  +6.38%  '
$ws.Range('D44').Value = '42.65'
$ws.Range('E44').Value = '  +0.51%  '
$ws.Range('D45').Value = '0.999'
$ws.Range('E45').Value = '  -0.09%  '
$ws.Range('D46').Value = '4.45'
$ws.Range('E46').Value = '  +2.20%  '
$ws.Range('E47').Value = '  +4.52%  '
$ws.Range('E48').Value = '  +1.18%  '
$ws.Range('D49').Value = '2.481.33'
$ws.Range('E49').Value = '  +11.91%  '
$ws.Range('E50').Value = '  +2.81%  '
$ws.Range('D51').Value = '0.866'
$ws.Range('E51').Value = '  +7.73%  '
